$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = "Chitradurga"
$ws.Range("G10").Value = "Chitradurga"
$ws.Range("G12").Value = "Chitradurga"
$ws.Range("G24").Value = "Chitradurga"
$ws.Range("G52").Value = "Chitradurga"
